$d = $word.ActiveDocument

# 1. Split "...the tests also create another " into
#    "...the tests " + "I included " + "also create another "
$d.Content.Find.Execute("the tests also create another", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the tests I included also create another", 2)

# 2. Add <w:semiHidden/> to the DefaultParagraphFont character style
$style = $d.Styles("Default Paragraph Font")
$style.SemiHidden = $true
